$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column F (description) to fit the expanded project list
$ws.Columns.Item(6).ColumnWidth = 78

# Append repository links (column C) for all newly catalogued projects, rows 60-102
$ws.Range('C60').Value = 'https://github.com/jobovy/apogee'
$ws.Range('C61').Value = 'https://github.com/qedsoftware/brukeropusreader'
$ws.Range('C62').Value = 'https://github.com/BIOS-IMASL/cheshift'
$ws.Range('C63').Value = 'https://github.com/chianti-atomic/ChiantiPy'
$ws.Range('C64').Value = 'https://github.com/PINT-NMR/COMPASS'
$ws.Range('C65').Value = 'https://github.com/CCampJr/CRIkit2'
$ws.Range('C66').Value = 'https://github.com/mretegan/crispy'
$ws.Range('C67').Value = 'https://github.com/prajwel/curvit'
$ws.Range('C68').Value = 'https://github.com/Kongsea/DeepSpectroscopy'
$ws.Range('C69').Value = 'https://github.com/rcthomas/es'
$ws.Range('C70').Value = 'https://github.com/Farseer-NMR/FarSeer-NMR'
$ws.Range('C71').Value = 'https://github.com/nasa-jpl/LiveViewOpenSource'
$ws.Range('C72').Value = 'https://github.com/CMET-UGent/MicroRaman'
$ws.Range('C73').Value = 'https://github.com/DeepanshS/mrsimulator'
$ws.Range('C74').Value = 'https://github.com/pnnl/nmrfit'
$ws.Range('C75').Value = 'https://github.com/nsherry4/Peakaboo'
$ws.Range('C76').Value = 'https://github.com/TheAstroFactory/pydis'
$ws.Range('C77').Value = 'https://github.com/sametz/pydnmr'
$ws.Range('C78').Value = 'https://github.com/L-Siemons/PyRoShift'
$ws.Range('C79').Value = 'https://github.com/ibcp/pyspectra'
$ws.Range('C80').Value = 'https://github.com/pycroscopy/pyUSID'
$ws.Range('C81').Value = 'https://github.com/radis/radis'
$ws.Range('C82').Value = 'https://github.com/DerekKaknes/raman'
$ws.Range('C83').Value = 'https://github.com/raman-noodles/Raman-noodles'
$ws.Range('C84').Value = 'https://github.com/charlesll/rampy'
$ws.Range('C85').Value = 'https://github.com/LlucSF/rCRSI'
$ws.Range('C86').Value = 'https://github.com/paris-saclay-cds/specio'
$ws.Range('C87').Value = 'https://github.com/astropy/specreduce'
$ws.Range('C88').Value = 'https://github.com/cheminfo-js/spectra-data'
$ws.Range('C89').Value = 'https://github.com/clerk67/spectra-formatter'
$ws.Range('C90').Value = 'https://github.com/charlesll/Spectra.jl'
$ws.Range('C91').Value = 'https://github.com/PlasmaPy/SpectroscoPyx'
$ws.Range('C92').Value = 'https://github.com/spacetelescope/specviz'
$ws.Range('C93').Value = 'https://github.com/workflow4metabolomics/tools-metabolomics'
$ws.Range('C94').Value = 'https://github.com/VespucciProject/Vespucci'
$ws.Range('C95').Value = 'https://github.com/archaeological-research-facility/web_geochemistry'
$ws.Range('C96').Value = 'https://github.com/megbedell/wobble'
$ws.Range('C97').Value = 'https://github.com/workflow4metabolomics/workflow4metabolomics'
$ws.Range('C98').Value = 'https://github.com/wcchu/XL-e'
$ws.Range('C99').Value = 'https://github.com/brianlevay/xrf_filetools'
$ws.Range('C100').Value = 'https://github.com/dylarm/xrf_parser'
$ws.Range('C101').Value = 'https://github.com/wojdyr/xylib'
$ws.Range('C102').Value = 'https://github.com/yokochi47/xyza2pipe'

# Fill in the remaining metadata for the fully documented new projects, rows 60-75
# apogee
$ws.Range('A60').Value = 'apogee'
$ws.Range('B60').Value = 'https://github.com/jobovy/apogee'
$ws.Range('D60').Value = 'https://arxiv.org/abs/1510.06745'
$ws.Range('E60').Value = 'Py'
$ws.Range('F60').Value = 'Tools for working with APOGEE data'
$ws.Range('G60').Value = 'IR'
$ws.Range('H60').Value = 'bovy at astro dot utoronto dot ca'

# brukeropusreader
$ws.Range('A61').Value = 'brukeropusreader'
$ws.Range('B61').Value = 'https://github.com/qedsoftware/brukeropusreader'
$ws.Range('E61').Value = 'Py'
$ws.Range('F61').Value = 'Read Bruker OPUS files'
$ws.Range('G61').Value = 'Data Sharing'
$ws.Range('H61').Value = 'brukeropusreader-dev@qed.ai'

# CheShift
$ws.Range('A62').Value = 'CheShift'
$ws.Range('B62').Value = 'https://github.com/BIOS-IMASL/cheshift'
$ws.Range('D62').Value = 'https://www.pnas.org/content/110/42/16826'
$ws.Range('E62').Value = 'Py'
$ws.Range('F62').Value = 'Prediction 13C shifts in proteins'
$ws.Range('G62').Value = 'NMR (13C)'

# ChiantiPy
$ws.Range('A63').Value = 'ChiantiPy'
$ws.Range('B63').Value = 'https://github.com/chianti-atomic/ChiantiPy'
$ws.Range('E63').Value = 'Py'
$ws.Range('F63').Value = 'Calculate radiative properties of astrophysical plasmas'
$ws.Range('G63').Value = 'UV-Vis-IR'

# COMPASS
$ws.Range('A64').Value = 'COMPASS'
$ws.Range('B64').Value = 'https://github.com/PINT-NMR/COMPASS'
$ws.Range('D64').Value = 'https://journals.plos.org/ploscompbiol/article?id=10.1371/journal.pcbi.1004022'
$ws.Range('E64').Value = 'Qt'
$ws.Range('F64').Value = 'Protein backbone assignments from triple-resonance peak lists'
$ws.Range('G64').Value = 'NMR'

# CRIkit2
$ws.Range('A65').Value = 'CRIkit2'
$ws.Range('B65').Value = 'https://github.com/CCampJr/CRIkit2'
$ws.Range('D65').Value = 'https://www.ncbi.nlm.nih.gov/pubmed/28819335'
$ws.Range('E65').Value = 'Py'
$ws.Range('F65').Value = 'Hyperspectral imaging toolkit'
$ws.Range('G65').Value = 'Raman'
$ws.Range('H65').Value = 'charles.camp@nist.gov'

# crispy
$ws.Range('A66').Value = 'crispy'
$ws.Range('B66').Value = 'http://www.esrf.eu/computing/scientific/crispy/'
$ws.Range('E66').Value = 'Py'
$ws.Range('F66').Value = 'Calculate & plot core-level spectra'

# curvit
$ws.Range('A67').Value = 'curvit'
$ws.Range('B67').Value = 'https://github.com/prajwel/curvit'
$ws.Range('E67').Value = 'Py'
$ws.Range('F67').Value = 'Create light curves from UV imaging telescope data'
$ws.Range('G67').Value = 'UV'

# DeepSpectroscopy
$ws.Range('A68').Value = 'DeepSpectroscopy'
$ws.Range('B68').Value = 'https://github.com/Kongsea/DeepSpectroscopy'
$ws.Range('E68').Value = 'Py'
$ws.Range('F68').Value = 'Spectroscopy with deep learning'

# es
$ws.Range('A69').Value = 'es'
$ws.Range('B69').Value = 'https://c3.lbl.gov/es/'
$ws.Range('E69').Value = 'Py'
$ws.Range('F69').Value = 'Elementary supernova spectrum synthesis'
$ws.Range('G69').Value = 'UV-Vis'

# Farseer-NMR
$ws.Range('A70').Value = 'Farseer-NMR'
$ws.Range('B70').Value = 'https://farseer-nmr.github.io/FarSeer-NMR/'
$ws.Range('D70').Value = 'https://link.springer.com/article/10.1007/s10858-018-0182-5'
$ws.Range('E70').Value = 'Py'
$ws.Range('F70').Value = 'Analysis & plotting of biological NMR peak lists'
$ws.Range('G70').Value = 'NMR (2D)'

# LiveView
$ws.Range('A71').Value = 'LiveView'
$ws.Range('B71').Value = 'https://github.com/nasa-jpl/LiveViewOpenSource'
$ws.Range('E71').Value = 'C++'
$ws.Range('F71').Value = 'Tools for imaging spectrometer development'
$ws.Range('G71').Value = 'UV-VIS-IR'
$ws.Range('H71').Value = 'Jacqueline.Ryan@jpl.nasa.gov'

# MicroRaman
$ws.Range('A72').Value = 'MicroRaman'
$ws.Range('B72').Value = 'https://www.ncbi.nlm.nih.gov/pubmed/29909167'
$ws.Range('D72').Value = 'https://www.ncbi.nlm.nih.gov/pubmed/29909167'
$ws.Range('E72').Value = 'R'
$ws.Range('F72').Value = 'Process microbial Raman spectra'
$ws.Range('G72').Value = 'Raman'
$ws.Range('H72').Value = 'frederiekmaarten.kerckhof@ugent.be'

# mrsimulator
$ws.Range('A73').Value = 'mrsimulator'
$ws.Range('B73').Value = 'https://deepanshs.github.io/mrsimulator/'
$ws.Range('E73').Value = 'Py'
$ws.Range('F73').Value = 'Toolbox for simulating NMR spectra'
$ws.Range('G73').Value = 'NMR'

# nmrfit
$ws.Range('A74').Value = 'nmrfit'
$ws.Range('B74').Value = 'https://github.com/pnnl/nmrfit'
$ws.Range('E74').Value = 'Py'
$ws.Range('F74').Value = 'qNMR through least-squares fitting'
$ws.Range('G74').Value = 'NMR'

# Peakaboo
$ws.Range('A75').Value = 'Peakaboo'
$ws.Range('B75').Value = 'https://peakaboo.org/'
$ws.Range('E75').Value = 'Java'
$ws.Range('F75').Value = 'XRF analysis'
$ws.Range('G75').Value = 'XRF'

# Leave the cursor where the author left off editing
$ws.Range("C76").Select()
